$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.036329779833898
$ws.Cells.Item(2, 4).Value = 1.043767832648233
$ws.Cells.Item(2, 5).Value = 1.035267632337253
$ws.Cells.Item(2, 6).Value = 1.052771044941338
$ws.Cells.Item(2, 9).Value = 1.040289004610983
$ws.Cells.Item(2, 10).Value = 1.04143868769633
$ws.Cells.Item(2, 11).Value = 1.046540769436517
$ws.Cells.Item(2, 12).Value = 1.038064725921372
$ws.Cells.Item(2, 13).Value = 1.055518856626008
$ws.Cells.Item(2, 14).Value = 1.005712725503983
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.037154938967109
$ws.Cells.Item(3, 4).Value = 1.044400786081595
$ws.Cells.Item(3, 5).Value = 1.035966107966737
$ws.Cells.Item(3, 6).Value = 1.053541225477026
$ws.Cells.Item(3, 9).Value = 1.040476963330802
$ws.Cells.Item(3, 10).Value = 1.041908392460722
$ws.Cells.Item(3, 11).Value = 1.046985290444099
$ws.Cells.Item(3, 12).Value = 1.038572856107397
$ws.Cells.Item(3, 13).Value = 1.056102058745903
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.037689606571
$ws.Cells.Item(4, 4).Value = 1.044810997294741
$ws.Cells.Item(4, 5).Value = 1.036419076330318
$ws.Cells.Item(4, 6).Value = 1.054040581820184
$ws.Cells.Item(4, 9).Value = 1.040597794981727
$ws.Cells.Item(4, 10).Value = 1.042212384957955
$ws.Cells.Item(4, 11).Value = 1.04727288220906
$ws.Cells.Item(4, 12).Value = 1.038901981688915
$ws.Cells.Item(4, 13).Value = 1.0564797638206
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.037914554885054
$ws.Cells.Item(5, 4).Value = 1.044983603355491
$ws.Cells.Item(5, 5).Value = 1.036609743612425
$ws.Cells.Item(5, 6).Value = 1.054250748094392
$ws.Cells.Item(5, 9).Value = 1.040648402667038
$ws.Cells.Item(5, 10).Value = 1.042340196920781
$ws.Cells.Item(5, 11).Value = 1.047393774112006
$ws.Cells.Item(5, 12).Value = 1.039040423958177
$ws.Cells.Item(5, 13).Value = 1.056638629389677
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.0379523348517
$ws.Cells.Item(6, 4).Value = 1.045012593605821
$ws.Cells.Item(6, 5).Value = 1.036641771469762
$ws.Cells.Item(6, 6).Value = 1.054286049766595
$ws.Cells.Item(6, 9).Value = 1.040656888759721
$ws.Cells.Item(6, 10).Value = 1.042361657865095
$ws.Cells.Item(6, 11).Value = 1.047414071670041
$ws.Cells.Item(6, 12).Value = 1.039063673548612
$ws.Cells.Item(6, 13).Value = 1.056665308163795
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.037692611659696
$ws.Cells.Item(7, 4).Value = 1.044813303065491
$ws.Cells.Item(7, 5).Value = 1.036421623097678
$ws.Cells.Item(7, 6).Value = 1.054043389144092
$ws.Cells.Item(7, 9).Value = 1.040598471951383
$ws.Cells.Item(7, 10).Value = 1.04221409273624
$ws.Cells.Item(7, 11).Value = 1.047274497620176
$ws.Cells.Item(7, 12).Value = 1.038903831257127
$ws.Cells.Item(7, 13).Value = 1.056481886284332
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.036608493078939
$ws.Cells.Item(8, 4).Value = 1.043981606924517
$ws.Cells.Item(8, 5).Value = 1.035503475596544
$ws.Cells.Item(8, 6).Value = 1.053031123168483
$ws.Cells.Item(8, 9).Value = 1.040352689299129
$ws.Cells.Item(8, 10).Value = 1.041597412948472
$ws.Cells.Item(8, 11).Value = 1.046691005483519
$ws.Cells.Item(8, 12).Value = 1.03823638140808
$ws.Cells.Item(8, 13).Value = 1.05571588228415
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.034703841611125
$ws.Cells.Item(9, 4).Value = 1.042521097776186
$ws.Cells.Item(9, 5).Value = 1.033893386682368
$ws.Cells.Item(9, 6).Value = 1.051255116269124
$ws.Cells.Item(9, 9).Value = 1.039913570626623
$ws.Cells.Item(9, 10).Value = 1.040511279828679
$ws.Cells.Item(9, 11).Value = 1.045662548614262
$ws.Cells.Item(9, 12).Value = 1.037062851852264
$ws.Cells.Item(9, 13).Value = 1.054368719974086
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.033438015363137
$ws.Cells.Item(10, 4).Value = 1.041550929263659
$ws.Cells.Item(10, 5).Value = 1.032825355892337
$ws.Cells.Item(10, 6).Value = 1.050076435780486
$ws.Cells.Item(10, 9).Value = 1.039616824245299
$ws.Cells.Item(10, 10).Value = 1.039787631877622
$ws.Cells.Item(10, 11).Value = 1.04497681155588
$ws.Cells.Item(10, 12).Value = 1.03628233557423
$ws.Cells.Item(10, 13).Value = 1.053472483469316
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.032890855920897
$ws.Cells.Item(11, 4).Value = 1.041131690445304
$ws.Cells.Item(11, 5).Value = 1.032364182826989
$ws.Cells.Item(11, 6).Value = 1.049567343092941
$ws.Cells.Item(11, 9).Value = 1.039487390564837
$ws.Cells.Item(11, 10).Value = 1.039474405965692
$ws.Cells.Item(11, 11).Value = 1.04467987364383
$ws.Cells.Item(11, 12).Value = 1.035944817502022
$ws.Cells.Item(11, 13).Value = 1.053084868361035
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.032687761332046
$ws.Cells.Item(12, 4).Value = 1.040976096195779
$ws.Cells.Item(12, 5).Value = 1.032193078604478
$ws.Cells.Item(12, 6).Value = 1.049378438400345
$ws.Cells.Item(12, 9).Value = 1.039439172597501
$ws.Cells.Item(12, 10).Value = 1.039358079088583
$ws.Cells.Item(12, 11).Value = 1.044569577693724
$ws.Cells.Item(12, 12).Value = 1.035819517459186
$ws.Cells.Item(12, 13).Value = 1.052940961866453
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.03273131922959
$ws.Cells.Item(13, 4).Value = 1.041009465789698
$ws.Cells.Item(13, 5).Value = 1.032229772140377
$ws.Cells.Item(13, 6).Value = 1.049418950239846
$ws.Cells.Item(13, 9).Value = 1.039449521862316
$ws.Cells.Item(13, 10).Value = 1.039383030708818
$ws.Cells.Item(13, 11).Value = 1.044593236526433
$ws.Cells.Item(13, 12).Value = 1.035846391588395
$ws.Cells.Item(13, 13).Value = 1.052971827056795
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.03287406509466
$ws.Cells.Item(14, 4).Value = 1.041118826318826
$ws.Cells.Item(14, 5).Value = 1.032350035289037
$ws.Cells.Item(14, 6).Value = 1.049551724181885
$ws.Cells.Item(14, 9).Value = 1.039483407715036
$ws.Cells.Item(14, 10).Value = 1.039464789950634
$ws.Cells.Item(14, 11).Value = 1.044670756537857
$ws.Cells.Item(14, 12).Value = 1.035934458741147
$ws.Cells.Item(14, 13).Value = 1.05307297154665
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.032962034773704
$ws.Cells.Item(15, 4).Value = 1.041186224196214
$ws.Cells.Item(15, 5).Value = 1.032424159413631
$ws.Cells.Item(15, 6).Value = 1.049633556492636
$ws.Cells.Item(15, 9).Value = 1.039504267308033
$ws.Cells.Item(15, 10).Value = 1.039515167099819
$ws.Cells.Item(15, 11).Value = 1.044718519218783
$ws.Cells.Item(15, 12).Value = 1.035988729034557
$ws.Cells.Item(15, 13).Value = 1.053135299466775
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.033474348738584
$ws.Cells.Item(16, 4).Value = 1.041578770848103
$ws.Cells.Item(16, 5).Value = 1.032855989805624
$ws.Cells.Item(16, 6).Value = 1.050110249814422
$ws.Cells.Item(16, 9).Value = 1.039625394581997
$ws.Cells.Item(16, 10).Value = 1.039808422262776
$ws.Cells.Item(16, 11).Value = 1.04499651827197
$ws.Cells.Item(16, 12).Value = 1.036304745158966
$ws.Cells.Item(16, 13).Value = 1.053498218071799
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.03379596582223
$ws.Cells.Item(17, 4).Value = 1.041825234070545
$ws.Cells.Item(17, 5).Value = 1.033127212723403
$ws.Cells.Item(17, 6).Value = 1.050409612140233
$ws.Cells.Item(17, 9).Value = 1.039701123254762
$ws.Cells.Item(17, 10).Value = 1.03999240615956
$ws.Cells.Item(17, 11).Value = 1.045170898247718
$ws.Cells.Item(17, 12).Value = 1.036503095426426
$ws.Cells.Item(17, 13).Value = 1.053725992020902
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.033983651327082
$ws.Cells.Item(18, 4).Value = 1.04196907383487
$ws.Cells.Item(18, 5).Value = 1.03328553692295
$ws.Cells.Item(18, 6).Value = 1.050584348839724
$ws.Cells.Item(18, 9).Value = 1.039745203708234
$ws.Cells.Item(18, 10).Value = 1.040099732187604
$ws.Cells.Item(18, 11).Value = 1.045272610151843
$ws.Cells.Item(18, 12).Value = 1.036618833160317
$ws.Cells.Item(18, 13).Value = 1.053858893096572
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.03404766273307
$ws.Cells.Item(19, 4).Value = 1.04201813331649
$ws.Cells.Item(19, 5).Value = 1.033339542433821
$ws.Cells.Item(19, 6).Value = 1.050643950458792
$ws.Cells.Item(19, 9).Value = 1.039760218579051
$ws.Cells.Item(19, 10).Value = 1.040136329470019
$ws.Cells.Item(19, 11).Value = 1.045307291072937
$ws.Cells.Item(19, 12).Value = 1.036658304048909
$ws.Cells.Item(19, 13).Value = 1.053904216404901
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.033761449854751
$ws.Cells.Item(20, 4).Value = 1.041798782410296
$ws.Cells.Item(20, 5).Value = 1.033098100183343
$ws.Cells.Item(20, 6).Value = 1.050377480583539
$ws.Cells.Item(20, 9).Value = 1.039693007671199
$ws.Cells.Item(20, 10).Value = 1.039972665246376
$ws.Cells.Item(20, 11).Value = 1.045152189018472
$ws.Cells.Item(20, 12).Value = 1.036481809822249
$ws.Cells.Item(20, 13).Value = 1.053701549428443
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.03283202598376
$ws.Cells.Item(21, 4).Value = 1.041086618772514
$ws.Cells.Item(21, 5).Value = 1.03231461536311
$ws.Cells.Item(21, 6).Value = 1.049512620165227
$ws.Cells.Item(21, 9).Value = 1.039473433049264
$ws.Cells.Item(21, 10).Value = 1.039440713351358
$ws.Cells.Item(21, 11).Value = 1.044647928815068
$ws.Cells.Item(21, 12).Value = 1.035908523237491
$ws.Cells.Item(21, 13).Value = 1.053043185041146
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.03224849891841
$ws.Cells.Item(22, 4).Value = 1.040639604667065
$ws.Cells.Item(22, 5).Value = 1.031823142263965
$ws.Cells.Item(22, 6).Value = 1.048969977351416
$ws.Cells.Item(22, 9).Value = 1.039334565057663
$ws.Cells.Item(22, 10).Value = 1.03910636569424
$ws.Cells.Item(22, 11).Value = 1.044330881170959
$ws.Cells.Item(22, 12).Value = 1.035548476301124
$ws.Cells.Item(22, 13).Value = 1.052629656944994
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.032557757559984
$ws.Cells.Item(23, 4).Value = 1.040876503431027
$ws.Cells.Item(23, 5).Value = 1.032083573215774
$ws.Cells.Item(23, 6).Value = 1.049257534822366
$ws.Cells.Item(23, 9).Value = 1.039408258385587
$ws.Cells.Item(23, 10).Value = 1.039283598672706
$ws.Cells.Item(23, 11).Value = 1.044498953586445
$ws.Cells.Item(23, 12).Value = 1.035739305476252
$ws.Cells.Item(23, 13).Value = 1.052848836410999
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.03377704585162
$ws.Cells.Item(24, 4).Value = 1.041810734524053
$ws.Cells.Item(24, 5).Value = 1.033111254502767
$ws.Cells.Item(24, 6).Value = 1.050391999068934
$ws.Cells.Item(24, 9).Value = 1.039696675034737
$ws.Cells.Item(24, 10).Value = 1.039981585280467
$ws.Cells.Item(24, 11).Value = 1.045160642917068
$ws.Cells.Item(24, 12).Value = 1.036491427736713
$ws.Cells.Item(24, 13).Value = 1.053712593846942
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.035195552771242
$ws.Cells.Item(25, 4).Value = 1.0428980649282
$ws.Cells.Item(25, 5).Value = 1.034308696590069
$ws.Cells.Item(25, 6).Value = 1.051713327248885
$ws.Cells.Item(25, 9).Value = 1.040027801816837
$ws.Cells.Item(25, 10).Value = 1.040791999496987
$ws.Cells.Item(25, 11).Value = 1.045928452833286
$ws.Cells.Item(25, 12).Value = 1.037365920003747
$ws.Cells.Item(25, 13).Value = 1.054716671167548
